$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values stay text, matching the source data (some
# look like plain numbers/dates and Excel would otherwise auto-convert them).
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.360.11'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.871.30'
$ws.Range('E3').Value = '  -0.53%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '235.51'
$ws.Range('E5').Value = '  -0.90%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4665'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2837'
$ws.Range('E8').Value = '  +1.14%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06551'
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '21.08'
$ws.Range('E10').Value = '  +7.79%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07934'
$ws.Range('E11').Value = '  +2.65%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '97.26'
$ws.Range('E12').Value = '  -0.94%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.875.78'
$ws.Range('E13').Value = '  -0.20%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.154'
$ws.Range('E14').Value = '  +0.86%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.6748'
$ws.Range('E15').Value = '  +1.54%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '281.74'
$ws.Range('E16').Value = '  -1.08%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '30.371.99'
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '5.546'
$ws.Range('E18').Value = '  +4.74%  '
$ws.Range('E19').Value = '  +0.19%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.68'
$ws.Range('E20').Value = '  +1.20%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '2.121.56'
$ws.Range('E21').Value = '  -0.35%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.000007288'
$ws.Range('E22').Value = '  +0.23%  '
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.198'
$ws.Range('E24').Value = '  +0.36%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.290'
$ws.Range('E25').Value = '  +0.50%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '164.45'
$ws.Range('E26').Value = '  -1.73%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.07'
$ws.Range('E27').Value = '  +0.57%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.934'
$ws.Range('E28').Value = '  -2.29%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.354'
$ws.Range('E29').Value = '  -1.29%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.09689'
$ws.Range('E30').Value = '  -1.61%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.439'
$ws.Range('E31').Value = '  -0.16%  '
$ws.Range('E32').Value = '  -1.00%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.116'
$ws.Range('E33').Value = '  -1.30%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.04712'
$ws.Range('E34').Value = '  +1.01%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7054'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.717'
$ws.Range('E37').Value = '  +0.43%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01859'
$ws.Range('E38').Value = '  -0.38%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.336'
$ws.Range('E39').Value = '  -5.23%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.545'
$ws.Range('E40').Value = '  +1.16%  '
$ws.Range('E41').Value = '  +2.27%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.947'
$ws.Range('E42').Value = '  -0.23%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.8485'
$ws.Range('E43').Value = '  -2.22%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.4194'
$ws.Range('E44').Value = '  +0.55%  '
$ws.Range('E45').Value = '  +0.21%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '103.92'
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '7.213'
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.269'
$ws.Range('E48').Value = '  -0.47%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '936.91'
$ws.Range('E49').Value = '  -4.99%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '34.17'
$ws.Range('E50').Value = '  +0.62%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.1133'
$ws.Range('E51').Value = '  -2.19%  '
